$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.162.70"
$ws.Range("E2").Value = "  +3.75%  "
$ws.Range("D3").Value = "2.428.00"
$ws.Range("E3").Value = "  +0.73%  "
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "316.60"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +3.34%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "102.60"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +5.55%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +7.33%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "35.47"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("E12").Value = "  -2.13%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "18.14"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").Value = "2.808.91"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "2.433.11"
$ws.Range("E16").Value = "  +0.36%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.840"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").Value = "45.085.91"
$ws.Range("E18").Value = "  +3.44%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "12.27"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "6.36"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "0.0₃0921"
$ws.Range("E21").Value = "  +2.24%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "68.85"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +0.75%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "243.83"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("E26").Value = "  -0.06%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "25.51"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("E29").Value = "  -11.93%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "49.23"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +3.02%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "33.00"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +1.96%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "20.31"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +10.47%  "
$ws.Range("E33").Value = "  +5.75%  "
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("E35").Value = "  +0.27%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.0763"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +1.67%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "1.87"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  +0.90%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "2.85"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -2.45%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "124.51"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("E42").Value = "  +1.02%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "20.78"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("E44").Value = "  +2.00%  "
$ws.Range("D45").Value = "1.934.65"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("E46").Value = "  -2.91%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "2.92"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +3.07%  "
$ws.Range("E48").Value = "  +15.71%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "9.17"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -1.48%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "76.50"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +5.96%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "53.89"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +2.32%  "
